$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 4773.971207788346
$ws.Range("G2").Value = 4777.112380590292
$ws.Range("F3").Value = 3987.239464740323
$ws.Range("G3").Value = 3990.564506662769
$ws.Range("F4").Value = 4055.418951690293
$ws.Range("G4").Value = 4059.002931093609
$ws.Range("F5").Value = 3933.168334034614
$ws.Range("G5").Value = 3935.715158287871
$ws.Range("F6").Value = 5414.516684202974
$ws.Range("G6").Value = 5418.458930659206
$ws.Range("F7").Value = 3597.759096711617
$ws.Range("G7").Value = 3600.962143980478
$ws.Range("F9").Value = 3891.078051059782
$ws.Range("G9").Value = 3893.720235691846
$ws.Range("F10").Value = 3657.680846070072
$ws.Range("G10").Value = 3660.380942276402
$ws.Range("F11").Value = 3440.132937663493
$ws.Range("G11").Value = 3442.832787887457
$ws.Range("F12").Value = 4267.345420107421
$ws.Range("G12").Value = 4269.187835723823
$ws.Range("F13").Value = 3789.952202197283
$ws.Range("G13").Value = 3792.815214598481
$ws.Range("F14").Value = 2012.556666321642
$ws.Range("G14").Value = 2012.594478498077
$ws.Range("F15").Value = 2589.229308755906
$ws.Range("G15").Value = 2590.23771505391
$ws.Range("F16").Value = 2021.262818792059
$ws.Range("G16").Value = 2022.625072707048
$ws.Range("F17").Value = 2766.519778841928
$ws.Range("G17").Value = 2768.099407204899
$ws.Range("F18").Value = 2409.801161042567
$ws.Range("G18").Value = 2411.15327430685
$ws.Range("F19").Value = 2409.678313636498
$ws.Range("G19").Value = 2410.327671679977
$ws.Range("F20").Value = 2618.927183173505
$ws.Range("G20").Value = 2626.31516506914
$ws.Range("F21").Value = 2968.776355421598
$ws.Range("G21").Value = 2973.164533132437
$ws.Range("F22").Value = 2797.233735770491
$ws.Range("G22").Value = 2802.413381403152
$ws.Range("F23").Value = 2697.203882007739
$ws.Range("G23").Value = 2703.102716753501
$ws.Range("F24").Value = 2888.137963153737
$ws.Range("G24").Value = 2893.309954908437
$ws.Range("F25").Value = 3105.748039102325
$ws.Range("G25").Value = 3110.528692290507
